$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-20 08:31:28"
$ws.Range("E3").Value = "2026-02-20 08:31:31"
$ws.Range("E4").Value = "2026-02-20 08:31:33"
$ws.Range("E5").Value = "2026-02-20 08:31:36"
$ws.Range("E6").Value = "2026-02-20 08:31:38"
$ws.Range("E7").Value = "2026-02-20 08:31:41"
$ws.Range("E8").Value = "2026-02-20 08:31:43"
$ws.Range("E9").Value = "2026-02-20 08:31:46"
$ws.Range("E10").Value = "2026-02-20 08:31:48"
$ws.Range("E11").Value = "2026-02-20 08:31:50"
$ws.Range("E12").Value = "2026-02-20 08:31:53"
$ws.Range("E13").Value = "2026-02-20 08:31:55"
$ws.Range("E14").Value = "2026-02-20 08:31:57"
$ws.Range("E15").Value = "2026-02-20 08:32:00"
$ws.Range("E16").Value = "2026-02-20 08:32:02"
$ws.Range("E17").Value = "2026-02-20 08:32:05"
$ws.Range("E18").Value = "2026-02-20 08:32:07"
$ws.Range("E19").Value = "2026-02-20 08:32:10"
$ws.Range("E20").Value = "2026-02-20 08:32:12"
$ws.Range("E21").Value = "2026-02-20 08:32:15"
$ws.Range("E22").Value = "2026-02-20 08:32:17"
$ws.Range("E23").Value = "2026-02-20 08:32:19"
$ws.Range("E24").Value = "2026-02-20 08:32:22"
$ws.Range("E25").Value = "2026-02-20 08:32:24"
$ws.Range("E26").Value = "2026-02-20 08:32:27"
$ws.Range("E27").Value = "2026-02-20 08:32:29"
$ws.Range("E28").Value = "2026-02-20 08:32:31"
$ws.Range("E29").Value = "2026-02-20 08:32:34"
$ws.Range("E30").Value = "2026-02-20 08:32:36"
$ws.Range("E31").Value = "2026-02-20 08:32:39"
$ws.Range("E32").Value = "2026-02-20 08:32:41"
$ws.Range("E33").Value = "2026-02-20 08:32:44"
$ws.Range("E34").Value = "2026-02-20 08:32:46"
$ws.Range("E35").Value = "2026-02-20 08:32:49"
$ws.Range("E36").Value = "2026-02-20 08:32:51"
$ws.Range("E37").Value = "2026-02-20 08:32:53"
$ws.Range("E38").Value = "2026-02-20 08:32:56"
$ws.Range("E39").Value = "2026-02-20 08:32:58"
$ws.Range("E40").Value = "2026-02-20 08:33:01"
$ws.Range("E41").Value = "2026-02-20 08:33:03"
$ws.Range("E42").Value = "2026-02-20 08:33:05"
$ws.Range("E43").Value = "2026-02-20 08:33:08"
$ws.Range("E44").Value = "2026-02-20 08:33:10"
$ws.Range("E45").Value = "2026-02-20 08:33:13"
$ws.Range("E46").Value = "2026-02-20 08:33:15"
